$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Cells.Item(2, 10).Value = 2511
$ws.Cells.Item(3, 10).Value = 2571
$ws.Cells.Item(4, 10).Value = 591
$ws.Cells.Item(5, 10).Value = 196
$ws.Cells.Item(6, 10).Value = 3210
$ws.Cells.Item(7, 10).Value = 9079

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Cells.Item(3, 10).Value = 107
$ws.Cells.Item(5, 10).Value = 11
$ws.Cells.Item(6, 10).Value = 88
$ws.Cells.Item(7, 10).Value = 305

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Cells.Item(2, 10).Value = 42
$ws.Cells.Item(6, 10).Value = 32
$ws.Cells.Item(7, 10).Value = 128

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Cells.Item(3, 10).Value = 133
$ws.Cells.Item(5, 10).Value = 8

$ws = $wb.Worksheets.Item('Gage Park')
$ws.Cells.Item(4, 10).Value = 8
$ws.Cells.Item(7, 10).Value = 68

$ws = $wb.Worksheets.Item('South Deering')
$ws.Cells.Item(6, 10).Value = 35
$ws.Cells.Item(7, 10).Value = 86

$ws = $wb.Worksheets.Item('New City')
$ws.Cells.Item(2, 10).Value = 69
$ws.Cells.Item(4, 10).Value = 12

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Cells.Item(4, 10).Value = 40
$ws.Cells.Item(5, 10).Value = 22
$ws.Cells.Item(7, 10).Value = 272
$ws.Cells.Item(8, 10).Value = 563
$ws.Cells.Item(9, 10).Value = 57
$ws.Cells.Item(11, 10).Value = 127
$ws.Cells.Item(12, 10).Value = 20
$ws.Cells.Item(19, 10).Value = 293
$ws.Cells.Item(23, 10).Value = 94
$ws.Cells.Item(29, 10).Value = 528
$ws.Cells.Item(31, 10).Value = 68
$ws.Cells.Item(33, 10).Value = 366
$ws.Cells.Item(34, 10).Value = 47
$ws.Cells.Item(36, 10).Value = 135
$ws.Cells.Item(37, 10).Value = 305
$ws.Cells.Item(42, 10).Value = 351
$ws.Cells.Item(43, 10).Value = 83
$ws.Cells.Item(44, 10).Value = 74
$ws.Cells.Item(46, 10).Value = 29
$ws.Cells.Item(47, 10).Value = 78
$ws.Cells.Item(48, 10).Value = 89
$ws.Cells.Item(50, 10).Value = 51
$ws.Cells.Item(51, 10).Value = 121
$ws.Cells.Item(54, 10).Value = 183
$ws.Cells.Item(55, 10).Value = 109
$ws.Cells.Item(63, 10).Value = 38
$ws.Cells.Item(64, 10).Value = 61
$ws.Cells.Item(66, 10).Value = 22
$ws.Cells.Item(76, 10).Value = 129
$ws.Cells.Item(79, 10).Value = 277
$ws.Cells.Item(80, 10).Value = 19
$ws.Cells.Item(81, 10).Value = 9
$ws.Cells.Item(83, 10).Value = 222
$ws.Cells.Item(84, 10).Value = 86
$ws.Cells.Item(85, 10).Value = 423
$ws.Cells.Item(86, 10).Value = 51
$ws.Cells.Item(88, 10).Value = 91
$ws.Cells.Item(90, 10).Value = 99
$ws.Cells.Item(94, 10).Value = 74
$ws.Cells.Item(95, 10).Value = 139
$ws.Cells.Item(97, 10).Value = 56
$ws.Cells.Item(99, 10).Value = 128
$ws.Cells.Item(101, 10).Value = 9079

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Cells.Item(6, 10).Value = 66
$ws.Cells.Item(7, 10).Value = 222

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Cells.Item(3, 10).Value = 39
$ws.Cells.Item(6, 10).Value = 40
$ws.Cells.Item(7, 10).Value = 139

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Cells.Item(3, 10).Value = 106
$ws.Cells.Item(6, 10).Value = 127
$ws.Cells.Item(7, 10).Value = 366

$ws = $wb.Worksheets.Item('Loop')
$ws.Cells.Item(2, 10).Value = 46
$ws.Cells.Item(3, 10).Value = 36
$ws.Cells.Item(7, 10).Value = 183

$ws = $wb.Worksheets.Item('Englewood')
$ws.Cells.Item(2, 10).Value = 155
$ws.Cells.Item(3, 10).Value = 174
$ws.Cells.Item(5, 10).Value = 23
$ws.Cells.Item(6, 10).Value = 147
$ws.Cells.Item(7, 10).Value = 528

$ws = $wb.Worksheets.Item('Chatham')
$ws.Cells.Item(3, 10).Value = 79
$ws.Cells.Item(7, 10).Value = 293

$ws = $wb.Worksheets.Item('Irving Park')
$ws.Cells.Item(2, 10).Value = 27
$ws.Cells.Item(7, 10).Value = 74

$ws = $wb.Worksheets.Item('Lake View')
$ws.Cells.Item(4, 10).Value = 16
$ws.Cells.Item(7, 10).Value = 89

$ws = $wb.Worksheets.Item('River North')
$ws.Cells.Item(2, 10).Value = 18
$ws.Cells.Item(6, 10).Value = 72
$ws.Cells.Item(7, 10).Value = 129

$ws = $wb.Worksheets.Item('South Shore')
$ws.Cells.Item(3, 10).Value = 160
$ws.Cells.Item(7, 10).Value = 423

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Cells.Item(2, 10).Value = 69
$ws.Cells.Item(3, 10).Value = 74
$ws.Cells.Item(7, 10).Value = 351

$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Cells.Item(3, 10).Value = 18
$ws.Cells.Item(6, 10).Value = 58
$ws.Cells.Item(7, 10).Value = 109

$ws = $wb.Worksheets.Item('Jefferson Park')
$ws.Cells.Item(3, 10).Value = 5
$ws.Cells.Item(7, 10).Value = 29

$ws = $wb.Worksheets.Item('Douglas')
$ws.Cells.Item(2, 10).Value = 27
$ws.Cells.Item(4, 10).Value = 10
$ws.Cells.Item(5, 10).Value = 3
$ws.Cells.Item(7, 10).Value = 94

$ws = $wb.Worksheets.Item('Roseland')
$ws.Cells.Item(2, 10).Value = 79
$ws.Cells.Item(3, 10).Value = 102
$ws.Cells.Item(5, 10).Value = 6
$ws.Cells.Item(6, 10).Value = 75
$ws.Cells.Item(7, 10).Value = 277

$ws = $wb.Worksheets.Item('Near South Side')
$ws.Cells.Item(3, 10).Value = 13
$ws.Cells.Item(7, 10).Value = 61

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Cells.Item(2, 10).Value = 49
$ws.Cells.Item(7, 10).Value = 135

$ws = $wb.Worksheets.Item('Garfield Ridge')
$ws.Cells.Item(6, 10).Value = 18
$ws.Cells.Item(7, 10).Value = 47

$ws = $wb.Worksheets.Item('West Loop')
$ws.Cells.Item(3, 10).Value = 12
$ws.Cells.Item(7, 10).Value = 74

$ws = $wb.Worksheets.Item('Kenwood')
$ws.Cells.Item(3, 10).Value = 21
$ws.Cells.Item(7, 10).Value = 78

$ws = $wb.Worksheets.Item('Lincoln Square')
$ws.Cells.Item(6, 10).Value = 13
$ws.Cells.Item(7, 10).Value = 51

$ws = $wb.Worksheets.Item('North Center')
$ws.Cells.Item(6, 10).Value = 12
$ws.Cells.Item(7, 10).Value = 22

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Cells.Item(2, 10).Value = 44
$ws.Cells.Item(7, 10).Value = 127

$ws = $wb.Worksheets.Item('Avalon Park')
$ws.Cells.Item(2, 10).Value = 15
$ws.Cells.Item(7, 10).Value = 57

$ws = $wb.Worksheets.Item('West Town')
$ws.Cells.Item(2, 10).Value = 13
$ws.Cells.Item(7, 10).Value = 56

$ws = $wb.Worksheets.Item('United Center')
$ws.Cells.Item(2, 10).Value = 19
$ws.Cells.Item(3, 10).Value = 27
$ws.Cells.Item(7, 10).Value = 91

$ws = $wb.Worksheets.Item('Austin')
$ws.Cells.Item(3, 10).Value = 184
$ws.Cells.Item(7, 10).Value = 563

$ws = $wb.Worksheets.Item('Armour Square')
$ws.Cells.Item(2, 10).Value = 7
$ws.Cells.Item(7, 10).Value = 22

$ws = $wb.Worksheets.Item('Streeterville')
$ws.Cells.Item(4, 10).Value = 25
$ws.Cells.Item(7, 10).Value = 51

$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Cells.Item(6, 10).Value = 37
$ws.Cells.Item(7, 10).Value = 99

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Cells.Item(6, 10).Value = 35
$ws.Cells.Item(7, 10).Value = 121

$ws = $wb.Worksheets.Item('Hyde Park')
$ws.Cells.Item(6, 10).Value = 50
$ws.Cells.Item(7, 10).Value = 83

$ws = $wb.Worksheets.Item('Rush & Division')
$ws.Cells.Item(6, 10).Value = 13
$ws.Cells.Item(7, 10).Value = 19

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Cells.Item(4, 10).Value = 6
$ws.Cells.Item(7, 10).Value = 272

$ws = $wb.Worksheets.Item('Archer Heights')
$ws.Cells.Item(3, 10).Value = 11
$ws.Cells.Item(7, 10).Value = 40

$ws = $wb.Worksheets.Item('Beverly')
$ws.Cells.Item(6, 10).Value = 14
$ws.Cells.Item(7, 10).Value = 20

$ws = $wb.Worksheets.Item('Sauganash,Forest Glen')
$ws.Cells.Item(3, 10).Value = 2
$ws.Cells.Item(7, 10).Value = 9
